# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" pair recorded in the two detail rows
# of the statement (rows 16 and 17) were swapped: period 2208 belongs
# with the 37333 amount, and period 2209 belongs with the 40000 amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: period 2209/40000 -> 2208/37333
$ws.Range("E16").Value = "2208"
$ws.Range("F16").Value = 37333

# Row 17: period 2208/37333 -> 2209/40000
$ws.Range("E17").Value = "2209"
$ws.Range("F17").Value = 40000
